$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Geo")

# Strip the "www/images/" prefix from the local image filenames in column H (rows 2-16)
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Value2
    if ($old -like "www/images/*") {
        $cell.Value2 = $old.Substring(11)
    }
}

# Add the new Byss-Calle / Alvkarleby row
$ws.Range("A17").Value = "Sweden"
$ws.Range("B17").Value = "Älvkarleby"
$ws.Range("C17").Value = 60.571396900000003
$ws.Range("D17").Value = 17.412947200000001
$ws.Range("E17").Value = "1800 CE"
$ws.Range("F17").Value = "Bysse-Calle was a prolific composer of nyckelharpa tunes, and played at a lot of weddings."
$ws.Range("G17").Value = "http://matscarlsson.eu/familj/pdf/byss-calle.pdf"
$ws.Range("H17").Value = "bysscalle.png"
$ws.Range("I17").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/byssecalle.png"
$ws.Range("J17").Value = "public domain, drawing by Olaf Thunman"

$ws.Hyperlinks.Add($ws.Range("I17"), "https://reportingnotes.com/wp-content/uploads/2019/05/byssecalle.png")

$ws.Range("F18").Select()
